$p = $ppt.ActivePresentation

# 1. Update the version/date text on the title slide.
$titleSlide = $p.Slides.Item(1)
$subtitle = $titleSlide.Shapes.Item(2).TextFrame.TextRange
$subtitle.Replace("Versie 3.1.0-dev, 04-04-2023", "Versie 3.1.0-dev, 02-06-2023", 0, 0, 0) | Out-Null

# 2. Insert a new measure slide "M35: Het project hanteert een agile
#    architectuuraanpak" right before the current M06 slide (slide 20),
#    pushing M06 and everything after it down by one position.
#    We do this by duplicating the M06 slide (identical layout/shape
#    structure to the new M35 slide) and rewriting the text of the
#    original (now-first) copy to the M35 content, leaving the
#    untouched duplicate to carry on as M06 one position later.
$m06Slide = $p.Slides.Item(20)
$m06Slide.Duplicate() | Out-Null

$newSlide = $p.Slides.Item(20)

$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$oldTitle = "M06: Het project meet kwaliteitsnormen geautomatiseerd en frequent"
$newTitle = "M35: Het project hanteert een agile architectuuraanpak"
$titleRange.Replace($oldTitle, $newTitle, 0, 0, 0) | Out-Null

$bodyRange = $newSlide.Shapes.Item(2).TextFrame.TextRange
$oldBody = "Het voldoen aan de kwaliteitsnormen die geautomatiseerd gemeten kunnen worden, wordt frequent en minimaal één keer per dag gemeten."
$newBody = "Tijdens de voorfase verwerkt het project de door de opdrachtgever opgestelde projectstartarchitectuur (PSA) in een eerste versie van het softwarearchitectuurdocument (SAD). Tijdens de realisatiefase werkt het project het SAD bij op basis van nieuwe inzichten."
$bodyRange.Replace($oldBody, $newBody, 0, 0, 0) | Out-Null
